$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110 (weekly price update), shifting the
# existing rows 110-122 down to 111-123.
$ws.Rows.Item(110).Insert()

$ws.Cells.Item(110, 1).Value = 11
$ws.Cells.Item(110, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(110, 3).Value = "Bíobío"
$ws.Cells.Item(110, 4).Value = 45077
$ws.Cells.Item(110, 5).Value = 8
$ws.Cells.Item(110, 6).Value = 100112012
$ws.Cells.Item(110, 7).Value = "Espinaca"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 40
$ws.Cells.Item(110, 11).Value = 5000
$ws.Cells.Item(110, 12).Value = 5500
$ws.Cells.Item(110, 13).Value = 5250
$ws.Cells.Item(110, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(110, 15).Value = "Región Metropolitana"
$ws.Cells.Item(110, 16).Value = 525
$ws.Cells.Item(110, 17).Value = 10
$ws.Cells.Item(110, 18).Value = "Hortaliza"

# Apply the same date number format used by the rest of column D.
$ws.Cells.Item(110, 4).NumberFormat = $ws.Cells.Item(109, 4).NumberFormat
